$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$base = "https://klasma.github.io/LoggingDetectiveFiles/Logging_OSTERSUND/tillsynsmail/"

$ws.Range("Y2").Value = 'HYPERLINK("' + $base + 'A 30683-2023.docx"; "A 30683-2023")'
$ws.Range("Y3").Value = 'HYPERLINK("' + $base + 'A 32699-2023.docx"; "A 32699-2023")'
$ws.Range("Y4").Value = 'HYPERLINK("' + $base + 'A 29992-2023.docx"; "A 29992-2023")'
